$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2022-01-17) is inserted as row 89; the
# previously existing rows 89-94 shift down to 90-95 unchanged.
$ws.Rows("89:89").Insert()

$ws.Range("A89").Value = 5
$ws.Range("B89").Value = "Macroferia Regional de Talca"
$ws.Range("C89").Value = "Maule"
$ws.Range("D89").Value = 44578
$ws.Range("E89").Value = 7
$ws.Range("F89").Value = 100112030
$ws.Range("G89").Value = "Poroto granado"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 350
$ws.Range("K89").Value = 22000
$ws.Range("L89").Value = 22000
$ws.Range("M89").Value = 22000
$ws.Range("N89").Value = "`$/saco 25 kilos"
$ws.Range("O89").Value = "Región del Maule"
$ws.Range("P89").Value = 880
$ws.Range("Q89").Value = 25
$ws.Range("R89").Value = "Hortaliza"
